$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 150 (the empty divider row just
# before the weekly-summary block). This shifts rows 150-155 down to
# 151-156 and copies formatting from the row above, matching the target
# workbook's row 150 styles (s="21"/"13"/"2"/"15"/"19").
$ws.Rows.Item(150).Insert()

# Populate the freshly inserted row 150 with the new "Tablet Mein Kochbuch"
# time-tracking entry.
$ws.Range("A150").Value = 22
$ws.Range("B150").Value = "Interface Design"
$ws.Range("C150").Value = "MockUps"
$ws.Range("D150").Value = "[FEATURE]"
$ws.Range("E150").Value = "Tablet Mein Kochbuch"
$ws.Range("F150").Value = 44492
$ws.Range("G150").Value = 44481
$ws.Range("J150").Value = 0.5
$ws.Range("K150").Value = 0.60416666666666663
$ws.Range("I150").Formula = "=ROUNDUP(((SUM(K150-J150)*24*60/60)/0.25),0)*0.25"

# Move the active selection to the new bottom-right corner cell, same as
# the saved workbook view in the target file.
[void]$ws.Range("I156").Select()
